$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename factor label in A24 ("Charlson_withoutage" -> "Charlson_withage")
$ws.Range("A24").Value = "Charlson_withage"

# Update boxOdds / boxCILow / boxCIHigh values for rows 2-30
$ws.Range("B2").Value = 1.41
$ws.Range("C2").Value = 0.57
$ws.Range("D2").Value = 3.52
$ws.Range("B3").Value = 0.97
$ws.Range("C3").Value = 0.93
$ws.Range("D3").Value = 1.02
$ws.Range("B4").Value = 0.99
$ws.Range("C4").Value = 0.96
$ws.Range("D4").Value = 1.03
$ws.Range("B5").Value = 0.96
$ws.Range("C5").Value = 0.89
$ws.Range("D5").Value = 1.03
$ws.Range("B6").Value = 0.69
$ws.Range("C6").Value = 0.16
$ws.Range("D6").Value = 2.93
$ws.Range("B7").Value = 2.12
$ws.Range("C7").Value = 0.88
$ws.Range("D7").Value = 5.1
$ws.Range("B8").Value = 1.83
$ws.Range("C8").Value = 0.76
$ws.Range("D8").Value = 4.41
$ws.Range("B9").Value = 2.98
$ws.Range("C9").Value = 0.32
$ws.Range("D9").Value = 27.75
$ws.Range("B10").Value = 1.78
$ws.Range("C10").Value = 0.7
$ws.Range("D10").Value = 1.16
$ws.Range("B11").Value = 1.55
$ws.Range("C11").Value = 0.57
$ws.Range("D11").Value = 4.19
$ws.Range("B12").Value = 0.85
$ws.Range("C12").Value = 0.11
$ws.Range("D12").Value = 0
$ws.Range("B13").Value = 0.94
$ws.Range("C13").Value = 0.4
$ws.Range("D13").Value = 1.46
$ws.Range("B14").Value = 1.31
$ws.Range("C14").Value = 0.56
$ws.Range("D14").Value = 1.18
$ws.Range("B15").Value = 1.16
$ws.Range("C15").Value = 0.35
$ws.Range("D15").Value = 3.85
$ws.Range("B16").Value = 0.88
$ws.Range("C16").Value = 0.22
$ws.Range("D16").Value = 3.49
$ws.Range("B17").Value = 0.72
$ws.Range("C17").Value = 0.49
$ws.Range("D17").Value = 1.08
$ws.Range("B18").Value = 0.98
$ws.Range("C18").Value = 0.64
$ws.Range("D18").Value = 1.51
$ws.Range("B19").Value = 0.79
$ws.Range("C19").Value = 0.62
$ws.Range("D19").Value = 1.01
$ws.Range("B20").Value = 1.24
$ws.Range("C20").Value = 0.84
$ws.Range("D20").Value = 1.85
$ws.Range("B21").Value = 0.94
$ws.Range("C21").Value = 0.74
$ws.Range("D21").Value = 1.18
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0
$ws.Range("B23").Value = 18
$ws.Range("C23").Value = 4.93
$ws.Range("D23").Value = 65.75
$ws.Range("B24").Value = 1.04
$ws.Range("C24").Value = 0.79
$ws.Range("D24").Value = 1.36
$ws.Range("B25").Value = 1.05
$ws.Range("C25").Value = 0.46
$ws.Range("D25").Value = 2.38
$ws.Range("B26").Value = 0.66
$ws.Range("C26").Value = 0.12
$ws.Range("D26").Value = 1.16
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 1.11
$ws.Range("B28").Value = 0.7
$ws.Range("C28").Value = 0.22
$ws.Range("D28").Value = 1.14
$ws.Range("B29").Value = 0.83
$ws.Range("C29").Value = 0.59
$ws.Range("D29").Value = 1.18
$ws.Range("B30").Value = 1.08
$ws.Range("C30").Value = 0.47
$ws.Range("D30").Value = 2.46
